$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "work progress" mini-table in rows 17-22, mirroring the layout
# --- already used for the "workplan" table in rows 11-16 (same column
# --- formatting for E:G). Copy the formats down first, then set the
# --- actual values/formulas for the new table.
$ws.Range("E11:G16").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 16: relabel the "Total hours:" cell as "Total estimated hours:"
# (this frees up "Total hours:" to be reused below, for the new table).
$ws.Range("E16").Value = "Total estimated hours:"

# Row 17: headers for the new table
$ws.Range("F17").Value = "Total work hours"
$ws.Range("G17").Value = "Estimated total cost"

# Rows 18-21: worked hours rate per role + computed total cost
$ws.Range("F18").Value = 11.7
$ws.Range("G18").Formula = "=D12*F18"

$ws.Range("F19").Value = 12.28
$ws.Range("G19").Formula = "=D13*F19"

$ws.Range("F20").Value = 14.65
$ws.Range("G20").Formula = "=D14*F20"

$ws.Range("F21").Value = 11.93
$ws.Range("G21").Formula = "=D15*F21"

# Row 22: totals for the new table
$ws.Range("E22").Value = "Total hours:"
$ws.Range("F22").Formula = "=SUM(F18:F21)"
$ws.Range("G22").Formula = "=SUM(G18:G21)"

# Widen the (until-now unused) H:I columns, as seen in the saved file.
$ws.Range("H1:I1").ColumnWidth = 10.57

# Move the saved selection to G8.
$ws.Range("G8").Select() | Out-Null
